# Book1.xlsx edit:
#  - Sheet1!A1 changes from 1200000.0 to 123.0 (C1's SUM(A:A) formula
#    recalculates automatically from 1522685.7 to 322808.7)
#  - a new worksheet "Sheet0" is added after Sheet1, with A1 = "1asd23"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A1").Value = 123.0

$ws0 = $wb.Worksheets.Add($null, $ws1)
$ws0.Name = "Sheet0"
$ws0.Range("A1").Value = "1asd23"
